# Updated symbol list on Sun Jan 15 11:25:33 UTC 2023 with GitHub Actions
# Refresh Price / Volume(1h) columns with latest scraped values, and
# fix the BOLO / CoinbaseStockToken rows (46/47) which were previously
# swapped relative to their correct rank order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) and Volume 1h % (E) columns are stored as literal text in this
# sheet (e.g. "296.40", "-2.35%"), so force text format before assigning
# the new values to avoid Excel auto-converting them to numbers/percentages.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "296.40"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-2.35%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "31.20"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.05%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.105"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-2.30%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07336"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.10%"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-1.39%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.673"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "11.76%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.739"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.12%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9186"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.12%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1674"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.62%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07082"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-5.56%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08073"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.97%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02992"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.26%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09906"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.20%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001505"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.19%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006174"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-4.50%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.451"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.56%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.231"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.26%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3272"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.98%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1332"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.56%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.554"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "4.19%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04623"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2.18%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-4.33%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001216"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.25%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004426"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.10%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.19%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001873"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "7.68%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01696"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "1.88%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04415"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-1.94%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007198"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.15%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1329"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-1.12%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-7.98%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-18.15%"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-1.54%"
$ws.Range("B46").Value = "CoinbaseStockToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.01021"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-21.37%"
$ws.Range("B47").Value = "BOLO"
$ws.Range("C47").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.868"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-1.28%"
